# Update timetable room/section codes across the six timetable sheets.
# Changes:
#   MINOR: Generative Ai [C101] -> MINOR: Generative Ai [C102]
#   MINOR: VLSI [C101]          -> MINOR: VLSI [C102]
#   CS262 (Lab) [L207] -> [L106]  (or [L106] -> [L107] on the "PreMid_Section_B/PostMid_*" sheets)
#   CS263 (Lab) [L207] -> [L107]  (or [L106] -> [L207] on the "PreMid_Section_B/PostMid_*" sheets)

$wb = $excel.ActiveWorkbook

# --- Regular_Section_A ---
$ws = $wb.Worksheets.Item("Regular_Section_A")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value  = "CS263 (Lab) [L107]"
$ws.Range("E7").Value  = "CS263 (Lab) [L107]"
$ws.Range("D8").Value  = "CS262 (Lab) [L106]"
$ws.Range("D9").Value  = "CS262 (Lab) [L106]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- Regular_Section_B ---
$ws = $wb.Worksheets.Item("Regular_Section_B")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value  = "CS263 (Lab) [L107]"
$ws.Range("E7").Value  = "CS263 (Lab) [L107]"
$ws.Range("D8").Value  = "CS262 (Lab) [L106]"
$ws.Range("D9").Value  = "CS262 (Lab) [L106]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PreMid_Section_A ---
$ws = $wb.Worksheets.Item("PreMid_Section_A")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value  = "CS262 (Lab) [L106]"
$ws.Range("C8").Value  = "CS263 (Lab) [L107]"
$ws.Range("B9").Value  = "CS262 (Lab) [L106]"
$ws.Range("C9").Value  = "CS263 (Lab) [L107]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PreMid_Section_B ---
$ws = $wb.Worksheets.Item("PreMid_Section_B")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value  = "CS262 (Lab) [L107]"
$ws.Range("C8").Value  = "CS263 (Lab) [L207]"
$ws.Range("B9").Value  = "CS262 (Lab) [L107]"
$ws.Range("C9").Value  = "CS263 (Lab) [L207]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PostMid_Section_A ---
$ws = $wb.Worksheets.Item("PostMid_Section_A")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value  = "CS262 (Lab) [L107]"
$ws.Range("C8").Value  = "CS263 (Lab) [L207]"
$ws.Range("B9").Value  = "CS262 (Lab) [L107]"
$ws.Range("C9").Value  = "CS263 (Lab) [L207]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PostMid_Section_B ---
$ws = $wb.Worksheets.Item("PostMid_Section_B")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value  = "CS262 (Lab) [L107]"
$ws.Range("C8").Value  = "CS263 (Lab) [L207]"
$ws.Range("B9").Value  = "CS262 (Lab) [L107]"
$ws.Range("C9").Value  = "CS263 (Lab) [L207]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"
